$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the H1 title
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Candyfinity - a candy-themed slot game with various exciting features and winning potential. Play for free today!</w:t></w:r></w:p><w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2 & 3. Near the end of the document: drop the duplicated bold "Play
#        Candyfinity..." paragraph and rewrite the italic meta-description
#        paragraph into the feature-image prompt (keeping the italic run).
# ---------------------------------------------------------------------
$playAgainText = "Play Candyfinity Online Slot Game for Free!"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $playAgainText -and $para.Style.NameLocal -ne "Heading 1") {
        $para.Range.Delete()
        break
    }
}

$oldReadText = "Read our review of Candyfinity - a candy-themed slot game with various exciting features and winning potential. Play for free today!"
$newImageText = 'Create an eye-catching feature image for Candyfinity in cartoon style featuring a happy Maya warrior with glasses. The warrior should be surrounded by various gummy candies, lollipops, sugary glazes, and all kinds of sweets that are the ingredients for this tasty and exciting slot game. Use vibrant colors such as red, pink, and violet to capture the explosion of colors and shapes in the game. The image should also include the game logo "Candyfinity" in bold and playful font. Let your creativity shine to grab the attention of online slot game players looking for a fun and exciting game to play.'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $oldReadText) {
        $full = $para.Range
        $textRange = $d.Range($full.Start, $full.End - 1)
        $textRange.Text = $newImageText
        break
    }
}
